# Apply "repull data, push all data, mean calculation" edit:
# Updates the dSF column (F) values for each start-row in the giolito_lucas
# sheet to reflect the newly re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = 1
    4  = -2
    5  = 9
    6  = -2
    7  = -9
    8  = -2
    9  = -1
    10 = -3
    11 = 2
    12 = -7
    13 = -2
    14 = 3
    15 = -5
    16 = -1
    17 = -2
    18 = 5
    19 = -2
    20 = -1
    21 = 2
    22 = 3
    23 = 3
    24 = -2
    25 = 1
    26 = -3
    27 = 1
    28 = 1
    30 = -3
    31 = 3
    32 = -1
    33 = -2
    36 = 0
    37 = -1
    38 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
